# Swap the data of row 8 and row 9 for the columns that differ between them:
# A, B, D, E, F, G, H, Q, R, Z, AB
# (all other columns already hold identical values in both rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $cell8 = $ws.Range($col + "8")
    $cell9 = $ws.Range($col + "9")

    $val8 = $cell8.Value2
    $val9 = $cell9.Value2

    $cell8.Value2 = $val9
    $cell9.Value2 = $val8
}
